## Applies the "added harvard case classification" update:
## - swaps the average_doctor / average_doctor_old column headers (BP1/BQ1)
## - refreshes every "*_old" comparison column (Ada_old, Avey_old, Buoy_old,
##   K health_old, WebMD_old, doctor_MA_old, doctor_NJ_old, doctor_TH_old,
##   average_doctor/average_doctor_old) for rows 4-13 with recomputed stats
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.187
$ws.Range("AJ4").Value = 0.063
$ws.Range("AK4").Value = 0.251
$ws.Range("AU4").Value = 0.143
$ws.Range("AV4").Value = 0.026
$ws.Range("AW4").Value = 0.162
$ws.Range("BA4").Value = 1.955
$ws.Range("BB4").Value = 0.169
$ws.Range("BC4").Value = 0.411
$ws.Range("BG4").Value = 0.73
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.6879999999999999
$ws.Range("BN4").Value = 0.093
$ws.Range("BO4").Value = 0.305
$ws.Range("BP4").Value = 0.652
$ws.Range("BQ4").Value = 0.658
$ws.Range("E4").Value = 0.378
$ws.Range("F4").Value = 0.08500000000000001
$ws.Range("G4").Value = 0.292
$ws.Range("N4").Value = 0.4
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.254
$ws.Range("W4").Value = 0.217
$ws.Range("X4").Value = 0.101
$ws.Range("Y4").Value = 0.318
$ws.Range("AI5").Value = 0.229
$ws.Range("AJ5").Value = 0.094
$ws.Range("AK5").Value = 0.307
$ws.Range("AU5").Value = 0.297
$ws.Range("AV5").Value = 0.099
$ws.Range("AW5").Value = 0.315
$ws.Range("BA5").Value = 1.369
$ws.Range("BB5").Value = 0.08599999999999999
$ws.Range("BC5").Value = 0.293
$ws.Range("BG5").Value = 0.413
$ws.Range("BH5").Value = 0.052
$ws.Range("BM5").Value = 0.57
$ws.Range("BN5").Value = 0.078
$ws.Range("BO5").Value = 0.279
$ws.Range("BP5").Value = 0.456
$ws.Range("BQ5").Value = 0.455
$ws.Range("E5").Value = 0.49
$ws.Range("F5").Value = 0.104
$ws.Range("G5").Value = 0.322
$ws.Range("N5").Value = 0.755
$ws.Range("O5").Value = 0.08400000000000001
$ws.Range("P5").Value = 0.289
$ws.Range("W5").Value = 0.225
$ws.Range("X5").Value = 0.113
$ws.Range("Y5").Value = 0.337
$ws.Range("AI6").Value = 0.206
$ws.Range("AU6").Value = 0.193
$ws.Range("BA6").Value = 1.6
$ws.Range("BG6").Value = 0.528
$ws.Range("BM6").Value = 0.623
$ws.Range("BP6").Value = 0.533
$ws.Range("BQ6").Value = 0.535
$ws.Range("E6").Value = 0.427
$ws.Range("N6").Value = 0.523
$ws.Range("W6").Value = 0.221
$ws.Range("AI7").Value = 0.219
$ws.Range("AU7").Value = 0.244
$ws.Range("BA7").Value = 1.451
$ws.Range("BG7").Value = 0.452
$ws.Range("BM7").Value = 0.59
$ws.Range("BQ7").Value = 0.484
$ws.Range("E7").Value = 0.463
$ws.Range("N7").Value = 0.641
$ws.Range("W7").Value = 0.223
$ws.Range("AI8").Value = 0.211
$ws.Range("AJ8").Value = 0.093
$ws.Range("AK8").Value = 0.305
$ws.Range("AU8").Value = 0.236
$ws.Range("AV8").Value = 0.076
$ws.Range("AW8").Value = 0.275
$ws.Range("BA8").Value = 1.738
$ws.Range("BB8").Value = 0.137
$ws.Range("BC8").Value = 0.371
$ws.Range("BG8").Value = 0.576
$ws.Range("BH8").Value = 0.105
$ws.Range("BI8").Value = 0.324
$ws.Range("BM8").Value = 0.703
$ws.Range("BN8").Value = 0.07199999999999999
$ws.Range("BO8").Value = 0.268
$ws.Range("BP8").Value = 0.579
$ws.Range("BQ8").Value = 0.589
$ws.Range("E8").Value = 0.514
$ws.Range("F8").Value = 0.132
$ws.Range("G8").Value = 0.363
$ws.Range("N8").Value = 0.759
$ws.Range("O8").Value = 0.066
$ws.Range("P8").Value = 0.257
$ws.Range("W8").Value = 0.217
$ws.Range("X8").Value = 0.105
$ws.Range("Y8").Value = 0.325
$ws.Range("AI9").Value = 0.12
$ws.Range("AJ9").Value = 0.106
$ws.Range("AK9").Value = 0.325
$ws.Range("BA9").Value = 1.72
$ws.Range("BB9").Value = 0.246
$ws.Range("BC9").Value = 0.496
$ws.Range("BG9").Value = 0.62
$ws.Range("BH9").Value = 0.236
$ws.Range("BI9").Value = 0.485
$ws.Range("BM9").Value = 0.66
$ws.Range("BN9").Value = 0.224
$ws.Range("BO9").Value = 0.474
$ws.Range("BP9").Value = 0.573
$ws.Range("BQ9").Value = 0.57
$ws.Range("E9").Value = 0.44
$ws.Range("F9").Value = 0.246
$ws.Range("G9").Value = 0.496
$ws.Range("N9").Value = 0.64
$ws.Range("O9").Value = 0.23
$ws.Range("P9").Value = 0.48
$ws.Range("W9").Value = 0.12
$ws.Range("X9").Value = 0.106
$ws.Range("Y9").Value = 0.325
$ws.Range("AI10").Value = 0.24
$ws.Range("AJ10").Value = 0.182
$ws.Range("AK10").Value = 0.427
$ws.Range("AU10").Value = 0.22
$ws.Range("AV10").Value = 0.172
$ws.Range("AW10").Value = 0.414
$ws.Range("BA10").Value = 2.02
$ws.Range("BG10").Value = 0.66
$ws.Range("BH10").Value = 0.224
$ws.Range("BI10").Value = 0.474
$ws.Range("BM10").Value = 0.86
$ws.Range("BN10").Value = 0.12
$ws.Range("BO10").Value = 0.347
$ws.Range("BP10").Value = 0.673
$ws.Range("BQ10").Value = 0.699
$ws.Range("E10").Value = 0.5600000000000001
$ws.Range("F10").Value = 0.246
$ws.Range("G10").Value = 0.496
$ws.Range("N10").Value = 0.84
$ws.Range("O10").Value = 0.134
$ws.Range("P10").Value = 0.367
$ws.Range("W10").Value = 0.26
$ws.Range("X10").Value = 0.192
$ws.Range("Y10").Value = 0.439
$ws.Range("AI11").Value = 0.24
$ws.Range("AJ11").Value = 0.182
$ws.Range("AK11").Value = 0.427
$ws.Range("AU11").Value = 0.34
$ws.Range("AV11").Value = 0.224
$ws.Range("AW11").Value = 0.474
$ws.Range("BA11").Value = 2.02
$ws.Range("BG11").Value = 0.66
$ws.Range("BH11").Value = 0.224
$ws.Range("BI11").Value = 0.474
$ws.Range("BM11").Value = 0.86
$ws.Range("BN11").Value = 0.12
$ws.Range("BO11").Value = 0.347
$ws.Range("BP11").Value = 0.673
$ws.Range("BQ11").Value = 0.699
$ws.Range("E11").Value = 0.58
$ws.Range("F11").Value = 0.244
$ws.Range("G11").Value = 0.494
$ws.Range("N11").Value = 0.88
$ws.Range("O11").Value = 0.106
$ws.Range("P11").Value = 0.325
$ws.Range("W11").Value = 0.26
$ws.Range("X11").Value = 0.192
$ws.Range("Y11").Value = 0.439
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.647
$ws.Range("AV12").Value = 1.758
$ws.Range("AW12").Value = 1.326
$ws.Range("BA12").Value = 3.587
$ws.Range("BB12").Value = 0.32
$ws.Range("BC12").Value = 0.5659999999999999
$ws.Range("BG12").Value = 1.061
$ws.Range("BH12").Value = 0.057
$ws.Range("BI12").Value = 0.239
$ws.Range("BM12").Value = 1.326
$ws.Range("BN12").Value = 0.406
$ws.Range("BO12").Value = 0.637
$ws.Range("BP12").Value = 1.196
$ws.Range("BQ12").Value = 1.259
$ws.Range("E12").Value = 1.414
$ws.Range("F12").Value = 0.656
$ws.Range("G12").Value = 0.8100000000000001
$ws.Range("N12").Value = 1.652
$ws.Range("O12").Value = 1.618
$ws.Range("P12").Value = 1.272
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI13").Value = 1.4
$ws.Range("AJ13").Value = 0.398
$ws.Range("AK13").Value = 0.631
$ws.Range("AU13").Value = 2.505
$ws.Range("AV13").Value = 1.307
$ws.Range("AW13").Value = 1.143
$ws.Range("BA13").Value = 2.517
$ws.Range("BB13").Value = 0.31
$ws.Range("BC13").Value = 0.556
$ws.Range("BG13").Value = 0.625
$ws.Range("BH13").Value = 0.08500000000000001
$ws.Range("BI13").Value = 0.292
$ws.Range("BM13").Value = 0.989
$ws.Range("BN13").Value = 0.368
$ws.Range("BO13").Value = 0.607
$ws.Range("BP13").Value = 0.839
$ws.Range("BQ13").Value = 0.785
$ws.Range("E13").Value = 1.728
$ws.Range("F13").Value = 0.92
$ws.Range("G13").Value = 0.959
$ws.Range("N13").Value = 2.348
$ws.Range("O13").Value = 1.164
$ws.Range("P13").Value = 1.079
$ws.Range("W13").Value = 1.09
$ws.Range("X13").Value = 0.186
$ws.Range("Y13").Value = 0.431
